$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B "gauge number" text cells: strip the " / NNNNN" (or "/NNNNN") suffix,
# keeping the leading digits (and the separating space, when the original used
# " / " rather than a bare "/").
$gaugeUpdates = @(
    @{Row = 2; Text = "00410 "}
    @{Row = 3; Text = "00410 "}
    @{Row = 4; Text = "00453 "}
    @{Row = 5; Text = "00453 "}
    @{Row = 6; Text = "00480 "}
    @{Row = 7; Text = "00480 "}
    @{Row = 8; Text = "00733 "}
    @{Row = 9; Text = "00733 "}
    @{Row = 10; Text = "00744 "}
    @{Row = 11; Text = "00744 "}
    @{Row = 12; Text = "00801"}
    @{Row = 13; Text = "00801"}
    @{Row = 14; Text = "00810 "}
    @{Row = 15; Text = "00810 "}
    @{Row = 16; Text = "00834 "}
    @{Row = 17; Text = "00834 "}
    @{Row = 18; Text = "00841 "}
    @{Row = 19; Text = "00841 "}
    @{Row = 20; Text = "00892 "}
    @{Row = 21; Text = "00892 "}
    @{Row = 22; Text = "00898"}
    @{Row = 23; Text = "00898"}
    @{Row = 24; Text = "00901 "}
    @{Row = 25; Text = "00901 "}
    @{Row = 26; Text = "00906 "}
    @{Row = 27; Text = "00906 "}
    @{Row = 28; Text = "00911 "}
    @{Row = 29; Text = "00911 "}
    @{Row = 30; Text = "00975 "}
    @{Row = 31; Text = "00975 "}
    @{Row = 32; Text = "00978 "}
    @{Row = 33; Text = "00978 "}
    @{Row = 34; Text = "00996"}
    @{Row = 35; Text = "00996"}
    @{Row = 36; Text = "00998"}
    @{Row = 37; Text = "00998"}
    @{Row = 38; Text = "01002"}
    @{Row = 39; Text = "01002"}
    @{Row = 40; Text = "01004"}
    @{Row = 41; Text = "01004"}
    @{Row = 42; Text = "01006"}
    @{Row = 43; Text = "01006"}
    @{Row = 44; Text = "01010"}
    @{Row = 45; Text = "01010"}
    @{Row = 46; Text = "01021 "}
    @{Row = 47; Text = "01021 "}
    @{Row = 48; Text = "01030 "}
    @{Row = 49; Text = "01030 "}
    @{Row = 50; Text = "01034 "}
    @{Row = 51; Text = "01034 "}
    @{Row = 52; Text = "01038 "}
    @{Row = 53; Text = "01038 "}
    @{Row = 54; Text = "01040"}
    @{Row = 55; Text = "01040"}
    @{Row = 56; Text = "01044"}
    @{Row = 57; Text = "01044"}
    @{Row = 58; Text = "01044(5)"}
    @{Row = 59; Text = "01044(5)"}
    @{Row = 60; Text = "01058"}
    @{Row = 61; Text = "01058"}
    @{Row = 62; Text = "01060"}
    @{Row = 63; Text = "01060"}
    @{Row = 64; Text = "01070"}
    @{Row = 65; Text = "01070"}
    @{Row = 66; Text = "01077"}
    @{Row = 67; Text = "01077"}
    @{Row = 68; Text = "01079"}
    @{Row = 69; Text = "01079"}
    @{Row = 70; Text = "01085"}
    @{Row = 71; Text = "01085"}
    @{Row = 72; Text = "01089"}
    @{Row = 73; Text = "01089"}
    @{Row = 74; Text = "01111"}
    @{Row = 75; Text = "01111"}
    @{Row = 76; Text = "01128"}
    @{Row = 77; Text = "01128"}
    @{Row = 78; Text = "01137"}
    @{Row = 79; Text = "01137"}
    @{Row = 80; Text = "01151"}
    @{Row = 81; Text = "01151"}
    @{Row = 82; Text = "01172"}
    @{Row = 83; Text = "01172"}
    @{Row = 84; Text = "01182"}
    @{Row = 85; Text = "01182"}
    @{Row = 86; Text = "01211"}
    @{Row = 87; Text = "01211"}
    @{Row = 88; Text = "01215"}
    @{Row = 89; Text = "01215"}
    @{Row = 90; Text = "01220"}
    @{Row = 91; Text = "01220"}
    @{Row = 92; Text = "01225"}
    @{Row = 93; Text = "01225"}
    @{Row = 94; Text = "01236"}
    @{Row = 95; Text = "01236"}
    @{Row = 96; Text = "01245"}
    @{Row = 97; Text = "01245"}
    @{Row = 98; Text = "01250"}
    @{Row = 99; Text = "01250"}
    @{Row = 100; Text = "01258"}
    @{Row = 101; Text = "01258"}
    @{Row = 102; Text = "01260"}
    @{Row = 103; Text = "01260"}
    @{Row = 104; Text = "01284"}
    @{Row = 105; Text = "01284"}
    @{Row = 106; Text = "01289"}
    @{Row = 107; Text = "01289"}
    @{Row = 108; Text = "01298"}
    @{Row = 109; Text = "01298"}
    @{Row = 110; Text = "01312"}
    @{Row = 111; Text = "01312"}
    @{Row = 485; Text = "424"}
    @{Row = 486; Text = "424"}
)

foreach ($u in $gaugeUpdates) {
    $cell = $ws.Cells.Item($u.Row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Text
}

# The update loop that produced this data paired up row-groups by position and
# had an off-by-one: the "01044(5)" group (rows 56-57) and the "01044/01045" group
# (rows 58-59) got swapped, so the A (index) and C (NAM number) columns for those
# four rows also need to swap to match.
$a56 = $ws.Cells.Item(56, 1).Value2
$a57 = $ws.Cells.Item(57, 1).Value2
$a58 = $ws.Cells.Item(58, 1).Value2
$a59 = $ws.Cells.Item(59, 1).Value2
$c56 = $ws.Cells.Item(56, 3).Value2
$c57 = $ws.Cells.Item(57, 3).Value2
$c58 = $ws.Cells.Item(58, 3).Value2
$c59 = $ws.Cells.Item(59, 3).Value2

$ws.Cells.Item(56, 1).Value = $a58
$ws.Cells.Item(57, 1).Value = $a59
$ws.Cells.Item(58, 1).Value = $a56
$ws.Cells.Item(59, 1).Value = $a57

$ws.Cells.Item(56, 3).Value = $c58
$ws.Cells.Item(57, 3).Value = $c59
$ws.Cells.Item(58, 3).Value = $c56
$ws.Cells.Item(59, 3).Value = $c57
